{"js": "// \"aggiunta tavola dei contenuti\" - Word regenerates the hidden _Toc\n// bookmarks (used by the Table of Contents field to locate each\n// heading) whenever a TOC is (re)built. Renumber the six heading\n// bookmarks to their new _Toc names, keeping each bookmark's range\n// (and therefore its w:id slot) untouched.\nconst doc = context.document;\n\nconst renames = [\n  [\"_Toc5728332\", \"_Toc5728479\"],\n  [\"_Toc5728333\", \"_Toc5728480\"],\n  [\"_Toc5728334\", \"_Toc5728481\"],\n  [\"_Toc5728335\", \"_Toc5728482\"],\n  [\"_Toc5728336\", \"_Toc5728483\"],\n  [\"_Toc5728337\", \"_Toc5728484\"],\n];\n\nfor (const [oldName, newName] of renames) {\n  const range = doc.getBookmarkRange(oldName);\n  doc.deleteBookmark(oldName);\n  range.insertBookmark(newName);\n}\n\nawait context.sync();\n", "ps1": "# \"aggiunta tavola dei contenuti\" - Word regenerates the hidden _Toc\n# bookmarks (used by the Table of Contents field to locate headings)\n# whenever a TOC is (re)built. Renumber the six heading bookmarks to\n# the new _Toc names, keeping their w:id and range untouched.\n$d = $word.ActiveDocument\n\n$renames = @{\n    \"_Toc5728332\" = \"_Toc5728479\"\n    \"_Toc5728333\" = \"_Toc5728480\"\n    \"_Toc5728334\" = \"_Toc5728481\"\n    \"_Toc5728335\" = \"_Toc5728482\"\n    \"_Toc5728336\" = \"_Toc5728483\"\n    \"_Toc5728337\" = \"_Toc5728484\"\n}\n\nforeach ($oldName in $renames.Keys) {\n    $newName = $renames[$oldName]\n    $bm = $d.Bookmarks($oldName)\n    $rng = $bm.Range\n    $bm.Delete()\n    $d.Bookmarks.Add($newName, $rng)\n}\n"}
